$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 772
$ws1.Range("F5").Value  = 256
$ws1.Range("F7").Value  = 1096
$ws1.Range("F9").Value  = 18
$ws1.Range("F11").Value = 97
$ws1.Range("F12").Value = 1083
$ws1.Range("F15").Value = 724
$ws1.Range("F20").Value = 617
$ws1.Range("F23").Value = 1902
$ws1.Range("F24").Value = 481
$ws1.Range("F27").Value = 256
$ws1.Range("F28").Value = 2535
$ws1.Range("F31").Value = 650
$ws1.Range("F35").Value = 872
$ws1.Range("F36").Value = 1582
$ws1.Range("F39").Value = 511
$ws1.Range("F40").Value = 108

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 120

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 772
$ws4.Range("F7").Value  = 256
$ws4.Range("F9").Value  = 1096
$ws4.Range("F11").Value = 18
$ws4.Range("F13").Value = 97
$ws4.Range("F14").Value = 1083
$ws4.Range("F16").Value = 724
$ws4.Range("F19").Value = 120
$ws4.Range("F20").Value = 120
$ws4.Range("F25").Value = 617
$ws4.Range("F28").Value = 1902
$ws4.Range("F29").Value = 481
$ws4.Range("F32").Value = 2535
$ws4.Range("F38").Value = 650
$ws4.Range("F42").Value = 872
$ws4.Range("F43").Value = 1582
$ws4.Range("F46").Value = 511
$ws4.Range("F47").Value = 108
